$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(27, 0, 12, 0, 0, 0, 15)
    3  = @(27, 0, 13, 0, 0, 0, 14)
    4  = @(27, 1, 14, 1, 1, 0, 10)
    5  = @(27, 0, 12, 2, 0, 1, 12)
    6  = @(27, 0, 12, 3, 2, 1, 9)
    7  = @(27, 0, 14, 1, 1, 1, 10)
    8  = @(27, 0, 13, 1, 2, 1, 10)
    9  = @(27, 0, 12, 2, 0, 1, 12)
    10 = @(27, 0, 14, 1, 0, 1, 11)
    11 = @(27, 0, 12, 1, 0, 0, 14)
    12 = @(27, 0, 13, 0, 0, 0, 14)
    13 = @(27, 0, 12, 1, 0, 1, 13)
    14 = @(27, 0, 13, 0, 2, 3, 9)
    15 = @(27, 0, 15, 1, 1, 3, 7)
    16 = @(27, 0, 15, 0, 0, 0, 12)
    17 = @(27, 0, 16, 1, 0, 1, 9)
    18 = @(27, 0, 12, 4, 0, 1, 10)
    19 = @(27, 0, 14, 2, 0, 0, 11)
    20 = @(27, 0, 14, 3, 0, 0, 10)
    21 = @(27, 0, 12, 2, 0, 0, 13)
    22 = @(27, 0, 14, 1, 0, 1, 11)
    23 = @(27, 0, 12, 4, 0, 1, 10)
    24 = @(27, 0, 12, 0, 0, 5, 10)
    25 = @(27, 0, 13, 2, 2, 0, 10)
    26 = @(27, 0, 16, 2, 0, 1, 8)
    27 = @(27, 0, 12, 3, 1, 1, 10)
    28 = @(27, 0, 13, 1, 1, 0, 12)
    29 = @(27, 0, 11, 1, 0, 2, 13)
    30 = @(27, 0, 13, 2, 0, 2, 10)
    31 = @(27, 0, 15, 2, 0, 2, 8)
    32 = @(27, 0, 13, 2, 0, 0, 12)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
